$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.716.05'
$ws.Range("E2").Value = '  -2.59%  '
$ws.Range("D3").Value = '3.808.21'
$ws.Range("E3").Value = '  +0.48%  '
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '601.85'
$ws.Range("E5").Value = '  -2.41%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '169.56'
$ws.Range("E6").Value = '  -4.82%  '
$ws.Range("D7").Value = '3.805.24'
$ws.Range("E7").Value = '  +0.45%  '
$ws.Range("E8").Value = '  +0.04%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.530'
$ws.Range("E9").Value = '  +0.29%  '
$ws.Range("E10").Value = '  -4.75%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.18'
$ws.Range("E11").Value = '  -5.47%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.467'
$ws.Range("E12").Value = '  -4.14%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '38.79'
$ws.Range("E13").Value = '  -3.03%  '
$ws.Range("E14").Value = '  -3.87%  '
$ws.Range("D15").Value = '4.436.95'
$ws.Range("E15").Value = '  +0.41%  '
$ws.Range("D16").Value = '3.806.70'
$ws.Range("E16").Value = '  +0.45%  '
$ws.Range("D17").Value = '67.782.47'
$ws.Range("E17").Value = '  -2.64%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '7.26'
$ws.Range("E18").Value = '  -4.44%  '
$ws.Range("E19").Value = '  -3.68%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '17.40'
$ws.Range("E20").Value = '  +5.77%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '495.32'
$ws.Range("E21").Value = '  -3.17%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '9.53'
$ws.Range("E22").Value = '  +1.06%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.744'
$ws.Range("E23").Value = '  +0.99%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '85.92'
$ws.Range("E24").Value = '  -0.63%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.39'
$ws.Range("E25").Value = '  -4.38%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.0000146'
$ws.Range("E26").Value = '  +6.44%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '12.38'
$ws.Range("E27").Value = '  -4.31%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.17'
$ws.Range("E28").Value = '  -4.15%  '
$ws.Range("E29").Value = '  +0.08%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.99'
$ws.Range("E30").Value = '  -0.87%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.44'
$ws.Range("E31").Value = '  -3.88%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '32.93'
$ws.Range("E32").Value = '  +5.51%  '
$ws.Range("E33").Value = '  -2.07%  '
$ws.Range("E34").Value = '  -4.25%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.00'
$ws.Range("E35").Value = '  -0.04%  '
$ws.Range("E36").Value = '  -3.10%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.85'
$ws.Range("E37").Value = '  -4.85%  '
$ws.Range("B38").Value = 'Bittensor'
$ws.Range("C38").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '463.90'
$ws.Range("E38").Value = '  +1.07%  '
$ws.Range("B39").Value = 'Kaspa'
$ws.Range("C39").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.133'
$ws.Range("E39").Value = '  -5.64%  '
$ws.Range("B40").Value = 'TheGraph'
$ws.Range("C40").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.330'
$ws.Range("E40").Value = '  -2.94%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '49.20'
$ws.Range("E41").Value = '  -1.26%  '
$ws.Range("E42").Value = '  -3.43%  '
$ws.Range("E43").Value = '  -3.86%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '8.45'
$ws.Range("E44").Value = '  -1.62%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '40.51'
$ws.Range("E46").Value = '  -8.68%  '
$ws.Range("D47").Value = '2.849.90'
$ws.Range("E47").Value = '  -3.85%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '140.17'
$ws.Range("E48").Value = '  +0.84%  '
$ws.Range("E49").Value = '  -2.03%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '24.72'
$ws.Range("E50").Value = '  +14.42%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '25.91'
$ws.Range("E51").Value = '  -5.24%  '
